# daily auto push: 2025-09-29 09:32 UTC
# Append a new data row (row 35) to the log sheet: date, weekday, hour,
# ranking - same shape as every existing row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count   # 34 (header + 33 data rows)
$newRow  = $lastRow + 1               # 35

# Column A/B of the previous row ("2025/09/29" / "月") already hold exactly
# the literal text the new row needs. Assigning a plain string straight to
# .Value runs it through Excel's normal text-entry parsing, which treats
# "2025/09/29" as a date and silently rewrites the cell into a date serial
# - not the plain text stored by every other row in this column. Copying
# the existing cells and pasting values-only reuses their already-correct
# text content (and default "no special format" styling) without
# re-parsing it, so the date string survives untouched.
$ws.Range("A" + $lastRow + ":B" + $lastRow).Copy()
$ws.Range("A" + $newRow + ":B" + $newRow).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$excel.CutCopyMode = 0

$ws.Cells.Item($newRow, 3).Value = 18
$ws.Cells.Item($newRow, 4).Value = 201
